# Applies the "Penalty Reward System" forecast-window shift:
# - "Forecast Comparison" sheet: each week's Week_Start_Date (col B) rolls
#   forward by one week, and MyForecast (col D) is replaced with new values.
# - "Summary" sheet: aggregate metrics recomputed from the new MyForecast
#   column.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# Keep the date-like strings stored as plain text (matching the workbook's
# existing inline-string convention) instead of being auto-converted to
# date serials.
$ws1.Range("B2:B17").NumberFormat = "@"

$newDates = @(
    "2025-01-12", "2025-01-19", "2025-01-26", "2025-02-02",
    "2025-02-09", "2025-02-16", "2025-02-23", "2025-03-02",
    "2025-03-09", "2025-03-16", "2025-03-23", "2025-03-30",
    "2025-04-06", "2025-04-13", "2025-04-20", "2025-04-27"
)

$newForecast = @(8, 7, 7, 6, 6, 5, 5, 5, 6, 5, 4, 8, 8, 8, 8, 8)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $newDates[$i]
    $ws1.Cells.Item($row, 4).Value = $newForecast[$i]
}

# --- Summary sheet -------------------------------------------------------
# Values in column B are stored as text, even the numeric-looking ones, so
# force text formatting on each cell before writing its new value.
$summaryUpdates = @{
    2  = "2022-12-25 to 2025-01-05"
    6  = "77"
    9  = "104"
    10 = "49"
    11 = "28"
    12 = "8"
    14 = "4"
    15 = "2025-03-23"
}

foreach ($row in $summaryUpdates.Keys) {
    $cell = $ws2.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryUpdates[$row]
}
